$d = $word.ActiveDocument

# 1) Expand on "søke fakta" -> "søke fakta rundt campusene" and
#    "generell forståelse." -> "generell forståelse for omfanget."
$d.Content.Find.Execute(
    "søke fakta og få en mer generell forståelse.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "søke fakta rundt campusene, og få en mer generell forståelse for omfanget.",
    2) | Out-Null

# 2) Drop the trailing space after "brukervennlighet." at the end of the
#    idé-/konsept paragraph.
$d.Content.Find.Execute(
    "brukervennlighet. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "brukervennlighet.",
    2) | Out-Null

# 3) The bulleted "Fargekode"/"Handlingsplan" paragraphs were set up with a
#    hanging indent (left indent 720 twips / 0.5" with a -720 twip
#    first-line indent). Remove the first-line (hanging) indent so the
#    paragraphs simply keep the 0.5" left indent.
foreach ($p in $d.Paragraphs) {
    if ($p.ParagraphFormat.LeftIndent -eq 36) {
        $p.ParagraphFormat.LeftIndent = 36
    }
}
